$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.805.21"
$ws.Range("E2").Value = "  +6.67%  "
$ws.Range("D3").Value = "'3.621.82"
$ws.Range("E3").Value = "  +4.14%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'419.48"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "'129.56"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.652"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("D8").Value = "'3.609.33"
$ws.Range("E8").Value = "  +3.96%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.761"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +24.75%  "
$ws.Range("D12").Value = "'0.0000428"
$ws.Range("E12").Value = "  +83.57%  "
$ws.Range("D13").Value = "'41.85"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "'9.79"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "'4.204.84"
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "'3.617.02"
$ws.Range("E17").Value = "  +4.57%  "
$ws.Range("D18").Value = "'20.02"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").Value = "'67.826.58"
$ws.Range("E20").Value = "  +6.85%  "
$ws.Range("D21").Value = "'12.37"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'458.51"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "'88.82"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "'13.36"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'3.04"
$ws.Range("E25").Value = "  -6.85%  "
$ws.Range("D26").Value = "'10.08"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").Value = "'35.44"
$ws.Range("E27").Value = "  +6.49%  "
$ws.Range("D28").Value = "'3.24"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").Value = "'4.95"
$ws.Range("E29").Value = "  +3.91%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'12.21"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.72"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("D32").Value = "'0.119"
$ws.Range("E32").Value = "  +5.09%  "
$ws.Range("D33").Value = "'7.20"
$ws.Range("E33").Value = "  -4.00%  "
$ws.Range("D34").Value = "'0.157"
$ws.Range("E34").Value = "  -7.11%  "
$ws.Range("D35").Value = "'39.89"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'56.00"
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("D38").Value = "'0.0₃0794"
$ws.Range("E38").Value = "  +22.56%  "
$ws.Range("D39").Value = "'0.0490"
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("E40").Value = "  +8.76%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'2.73"
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'147.82"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "'2.91"
$ws.Range("E44").Value = "  -4.83%  "
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.60"
$ws.Range("E46").Value = "  +11.24%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'4.25"
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.167"
$ws.Range("E48").Value = "  +19.57%  "
$ws.Range("D49").Value = "'0.300"
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").Value = "'2.63"
$ws.Range("E51").Value = "  +12.90%  "
